# Remove C8 and C10 from the capacitor designator list (row 3) and re-count
# the affected quantity, per the commit message:
#   "removed C8 and C10 from all variants as tests showed that without
#    these the SPI can be clocked a little faster"
#
# Row 3 held the 6-piece group "C2, C3, C4, C5, C8, C10" (22p / 0805 caps);
# dropping C8 and C10 leaves a 4-piece group "C2, C3, C4, C5".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Designator list text for that BOM line (column B).
$ws.Range("B3").Value = "C2, C3, C4, C5"

# Quantity for that BOM line (column A) drops from 6 to 4.
$ws.Range("A3").Value = 4

# Row got very slightly shorter (13.5pt -> 13.4pt) when the sheet was re-saved.
$ws.Rows.Item(3).RowHeight = 13.4

# Selection cursor moved from G1 to B3 in the re-saved sheet.
$ws.Range("B3").Select()
